$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (G and H) before the current "tag_ids" column,
# shifting tag_ids/doc_ids (and all their data) two columns to the right.
$ws.Range("G1:H1").EntireColumn.Insert()

# Match the width used by the neighbouring "phone" column for the two
# newly inserted columns.
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 13.5

# Header labels for the new columns.
$ws.Range("G1").Value2 = "start_date"
$ws.Range("H1").Value2 = "end_date"

# New validity start/end date values for a few institutions.
$ws.Range("G5").Value2 = 2012
$ws.Range("H5").Value2 = 2023
$ws.Range("G8").Value2 = "2010/10"
$ws.Range("H10").Value2 = "2021/04"

# Grow the table ("Tableau1") so it covers the two new columns, then make
# sure the last two columns keep their original header names (the insert
# operation leaves them with generic placeholder names).
$tbl = $ws.ListObjects.Item(1)
$lastRow = $tbl.Range.Rows.Count
$tbl.Resize($ws.Range("A1:J" + $lastRow))
$tbl.ListColumns.Item(9).Range.Cells.Item(1,1).Value2 = "tag_ids"
$tbl.ListColumns.Item(10).Range.Cells.Item(1,1).Value2 = "doc_ids"

# Restore the active-cell selection for the frozen bottom-right pane.
$ws.Range("H11").Select()
